$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Column C ("Förändrad") for all existing data rows (2..425) moves from
#    2023-09-20 (45189) to 2023-09-21 (45190).
for ($r = 2; $r -le 425; $r++) {
    $ws.Cells.Item($r, 3).Value = 45190
}

# 2) The two trailing rows (425 "A 44201-2023" / 426 "A 44205-2023") swap
#    their content, and the "changed" date for both becomes 45190.
#    Row 425 becomes what used to be row 426 ...
$ws.Range("A425").Value = "A 44205-2023"
$ws.Range("B425").Value = 45188
$ws.Range("C425").Value = 45190
$ws.Range("D425").Value = "DALARNAS LÄN"
$ws.Range("E425").Value = "MORA"
$ws.Range("G425").Value = 0.7

# ... and row 426 becomes what used to be row 425, now with an explicit
#     row height (customHeight) matching the other data rows.
$ws.Range("A426").Value = "A 44201-2023"
$ws.Range("B426").Value = 45188
$ws.Range("C426").Value = 45190
$ws.Range("D426").Value = "DALARNAS LÄN"
$ws.Range("E426").Value = "MORA"
$ws.Range("G426").Value = 0.6
$ws.Rows.Item(426).RowHeight = 15

# 3) A brand-new row 427 is appended.
$ws.Range("A427").Value = "A 44106-2023"
$ws.Range("B427").Value = 45188
$ws.Range("C427").Value = 45190
$ws.Range("D427").Value = "DALARNAS LÄN"
$ws.Range("E427").Value = "MORA"
$ws.Range("G427").Value = 10.1
$ws.Range("H427").Value = 0
$ws.Range("I427").Value = 0
$ws.Range("J427").Value = 0
$ws.Range("K427").Value = 0
$ws.Range("L427").Value = 0
$ws.Range("M427").Value = 0
$ws.Range("N427").Value = 0
$ws.Range("O427").Value = 0
$ws.Range("P427").Value = 0
$ws.Range("Q427").Value = 0
$ws.Range("R427").WrapText = $true

# Dates (columns B & C) keep the yyyy-mm-dd date format used elsewhere in
# the sheet; re-apply it explicitly for the two touched/newly written rows
# so the style index matches the other date cells.
$ws.Range("B425").NumberFormat = "YYYY-MM-DD"
$ws.Range("C425").NumberFormat = "YYYY-MM-DD"
$ws.Range("B426").NumberFormat = "YYYY-MM-DD"
$ws.Range("C426").NumberFormat = "YYYY-MM-DD"
$ws.Range("B427").NumberFormat = "YYYY-MM-DD"
$ws.Range("C427").NumberFormat = "YYYY-MM-DD"
